# Adds "green hydrogen" and "low carbon hydrogen" rows to the
# DACD-energyintensity sheet, mirroring the existing "hydrogen" row
# (same formulas/styling), and leaves that sheet as the active tab
# with the new rows selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DACD-energyintensity")
$ws.Activate()

# Insert "green hydrogen" as new row 12, copying the formatting/formulas
# from row 11 ("hydrogen") so every column keeps the same $B<row> pattern.
$ws.Rows("11").Copy()
$ws.Rows("12").Insert()
$ws.Range("A12").Value = "green hydrogen"

# Insert "low carbon hydrogen" as new row 13, same approach.
$ws.Rows("11").Copy()
$ws.Rows("13").Insert()
$ws.Range("A13").Value = "low carbon hydrogen"

# Leave the two new rows selected and this sheet as the active tab,
# matching the saved view state of the edited workbook.
$ws.Rows("12:13").Select()
